$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scc_p5_generation")

# Row 64: section header label (new dataset name)
$ws.Cells.Item(64, 1).Value = "artificial_regional_rare_celltype_diverse"

# Row 65: column headers (method names). A65 stays blank like the other
# section header rows (e.g. A2, A11, ...), so copy the blank/default
# formatting from one of those existing blank cells instead of writing a
# value, to reproduce the empty <c r="A65"/> cell.
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(65, 1))
$ws.Cells.Item(65, 2).Value = "spotlight"
$ws.Cells.Item(65, 3).Value = "music"
$ws.Cells.Item(65, 4).Value = "cell2location"
$ws.Cells.Item(65, 5).Value = "RCTD"
$ws.Cells.Item(65, 6).Value = "stereoscope"

# Row 66: corr
$ws.Cells.Item(66, 1).Value = "corr"
$ws.Cells.Item(66, 2).Value = 0.882041761053729
$ws.Cells.Item(66, 3).Value = 0.8826079562193381
$ws.Cells.Item(66, 4).Value = 0.9179191506358779
$ws.Cells.Item(66, 5).Value = 0.8853453459031966
$ws.Cells.Item(66, 6).Value = 0.8759874626340003

# Row 67: RMSE
$ws.Cells.Item(67, 1).Value = "RMSE"
$ws.Cells.Item(67, 2).Value = 2.444914051816823
$ws.Cells.Item(67, 3).Value = 2.548776716957037
$ws.Cells.Item(67, 4).Value = 2.0235906875439844
$ws.Cells.Item(67, 5).Value = 2.2829123450279694
$ws.Cells.Item(67, 6).Value = 2.4243985502573473

# Row 68: accuracy
$ws.Cells.Item(68, 1).Value = "accuracy"
$ws.Cells.Item(68, 2).Value = 0.91
$ws.Cells.Item(68, 3).Value = 0.72
$ws.Cells.Item(68, 4).Value = 0.83
$ws.Cells.Item(68, 5).Value = 0.82
$ws.Cells.Item(68, 6).Value = 0.78

# Row 69: sensitivity
$ws.Cells.Item(69, 1).Value = "sensitivity"
$ws.Cells.Item(69, 2).Value = 0.9
$ws.Cells.Item(69, 3).Value = 0.98
$ws.Cells.Item(69, 4).Value = 1.0
$ws.Cells.Item(69, 5).Value = 0.98
$ws.Cells.Item(69, 6).Value = 1.0

# Row 70: specificity
$ws.Cells.Item(70, 1).Value = "specificity"
$ws.Cells.Item(70, 2).Value = 0.91
$ws.Cells.Item(70, 3).Value = 0.67
$ws.Cells.Item(70, 4).Value = 0.8
$ws.Cells.Item(70, 5).Value = 0.79
$ws.Cells.Item(70, 6).Value = 0.73

# Row 71: precision
$ws.Cells.Item(71, 1).Value = "precision"
$ws.Cells.Item(71, 2).Value = 0.66
$ws.Cells.Item(71, 3).Value = 0.36
$ws.Cells.Item(71, 4).Value = 0.49
$ws.Cells.Item(71, 5).Value = 0.47
$ws.Cells.Item(71, 6).Value = 0.41

# Row 72: F1
$ws.Cells.Item(72, 1).Value = "F1"
$ws.Cells.Item(72, 2).Value = 0.76
$ws.Cells.Item(72, 3).Value = 0.53
$ws.Cells.Item(72, 4).Value = 0.66
$ws.Cells.Item(72, 5).Value = 0.64
$ws.Cells.Item(72, 6).Value = 0.58
